# Add working set of sequences
# For a specific set of rows in Sheet1, columns F through N all hold the
# same shared string value ("N/A"). In the "before" workbook many of
# those rows were only populated through column F; this script extends
# each of those rows so that columns G through N are populated with the
# same "N/A" value as column F, matching the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(6,7,14,18,19,20,27,29,32,33,34,39,40,44,46,50,51,55,62,68,75,78,80,84,88,89,91,93,99,101,106,109,115,116,117,118,121,122,127,128,129,132,133,135,137,138,140,142,143,145,151,153,154,155,156,159,163,170,173,175,180,186,191,192)

foreach ($r in $rows) {
    for ($col = 7; $col -le 14; $col++) {
        $ws.Cells.Item($r, $col).Value = "N/A"
    }
}
